# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" values are stored as plain text (some contain two
# "." separators, e.g. "27.532.72", which are not valid numbers). To keep
# ALL Price cells as text - matching the original workbook - we force the
# cell format to Text before writing, so Excel does not auto-convert
# numeric-looking strings (like "1.00") into real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.532.72"
$ws.Range("E2").Value = "  -3.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.661.13"
$ws.Range("E3").Value = "  -3.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.69"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.513"
$ws.Range("E6").Value = "  -2.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.48"
$ws.Range("E8").Value = "  +2.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.264"
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0620"
$ws.Range("E10").Value = "  -2.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0877"
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.895.25"
$ws.Range("E12").Value = "  -3.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.664.34"
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.14"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.567"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.94"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.535.51"
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "240.74"
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0731"
$ws.Range("E19").Value = "  -2.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.66"
$ws.Range("E20").Value = "  -3.01%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.46"
$ws.Range("E22").Value = "  -3.08%  "
$ws.Range("E23").Value = "  -2.38%  "
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.07"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.24"
$ws.Range("E26").Value = "  -2.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.29"
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.112"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("E30").Value = "  -2.99%  "
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.460.19"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("E34").Value = "  -4.37%  "
$ws.Range("E35").Value = "  -4.19%  "
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("E37").Value = "  -5.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.574"
$ws.Range("E38").Value = "  -4.57%  "
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.76"
$ws.Range("E42").Value = "  -4.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.48"
$ws.Range("E43").Value = "  -3.05%  "
$ws.Range("E44").Value = "  -2.81%  "
$ws.Range("E45").Value = "  -2.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.803.02"
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.72"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.81"
$ws.Range("E48").Value = "  -1.71%  "
$ws.Range("E49").Value = "  -5.91%  "
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.85"
$ws.Range("E51").Value = "  -3.39%  "
